$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, shifting existing rows 100-130 down to 101-131
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new weekly price record
$ws.Cells.Item(100, 1).Value = 10
$ws.Cells.Item(100, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(100, 3).Value = "La Araucanía"
$ws.Cells.Item(100, 4).Value = 45215
$ws.Cells.Item(100, 5).Value = 9
$ws.Cells.Item(100, 6).Value = 300000001
$ws.Cells.Item(100, 7).Value = "Rabanito"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 55
$ws.Cells.Item(100, 11).Value = 9000
$ws.Cells.Item(100, 12).Value = 9000
$ws.Cells.Item(100, 13).Value = 9000
$ws.Cells.Item(100, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(100, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(100, 16).Value = 750
$ws.Cells.Item(100, 17).Value = 12
$ws.Cells.Item(100, 18).Value = "Hortaliza"
